# Apply the numeric corrections described in the diff.
# All target cells hold plain numeric literals (no formulas), so we
# just overwrite each one with its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (EMPOLI)
$ws.Range("G3").Value = 5.881453154875717
$ws.Range("K3").Value = 6.047762760072973

# Row 4 (FIORENTINA)
$ws.Range("G4").Value = 6.210031347962382
$ws.Range("K4").Value = 6.498906752095594

# Row 10 (LAZIO)
$ws.Range("B10").Value = 6.166666666666667
$ws.Range("C10").Value = 5.802607076350093
$ws.Range("D10").Value = 5.981481481481482
$ws.Range("E10").Value = 6.028268551236749
$ws.Range("F10").Value = 4.563725490196078
$ws.Range("G10").Value = 5.793296089385475
$ws.Range("H10").Value = 6.1364522417154
$ws.Range("I10").Value = 6.674911660777386
$ws.Range("J10").Value = 4.952614379084967
$ws.Range("K10").Value = 5.914733379801572

# Row 18 (MILAN)
$ws.Range("B18").Value = 6.205882352941177
$ws.Range("C18").Value = 5.87378640776699
$ws.Range("D18").Value = 5.9521484375
$ws.Range("E18").Value = 5.964406779661017
$ws.Range("F18").Value = 5.102941176470588
$ws.Range("G18").Value = 5.885436893203884
$ws.Range("H18").Value = 6.1650390625
$ws.Range("I18").Value = 6.430508474576271
$ws.Range("J18").Value = 5.528546712802773
$ws.Range("K18").Value = 6.048741634461304

# Row 21 (ROMA)
$ws.Range("G21").Value = 5.966861598440546
$ws.Range("K21").Value = 6.156860230498272
